$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the
# 9c9bafea row (row 3) to reflect the newly generated handback report time.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-04 14:52:48"

# zh-cn sheet: update Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the 9c9bafea row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-04 14:52:43"
$wsZhCn.Range("K3").Value = "2016-09-04 14:53:00"

# de-de sheet: update Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the 9c9bafea row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-04 14:52:48"
$wsDeDe.Range("K3").Value = "2016-09-04 14:53:12"
